# Actualización automática 2025-06-01 08:00:06
# Apply the "roll forward one month" update to the VENTA MENSUAL sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Column widths (C, D, E, F) ---
# Note: Excel's ColumnWidth property stores an extra 5/6 (0.8333...) padding
# compared to the raw OOXML "width" attribute, so subtract it here to land
# on the exact target widths of 11, 11, 10, 11.
$ws.Columns.Item(3).ColumnWidth = 11 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 11 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 10 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 11 - (5/6)

# --- Header row: shift month labels forward by one ---
$ws.Range("C1").Value = "marzo"
$ws.Range("D1").Value = "abril"
$ws.Range("E1").Value = "mayo"
$ws.Range("F1").Value = "junio"

# --- Data rows: shift the sales values forward by one column ---
$ws.Range("C4").Value = 12.48
$ws.Range("D4").Value = 0

$ws.Range("D6").Value = 3.47
$ws.Range("E6").Value = 0

$ws.Range("C7").Value = 12.48
$ws.Range("D7").Value = 3.47
$ws.Range("E7").Value = 0
